$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row containing "US Core PMO ServiceRequest Profile" (row 43).
# Deleting the entire row shifts all subsequent rows up by one, which also
# naturally decrements the running index stored in column A for every row
# below (since it's a simple count), matching the target diff exactly.
$ws.Rows.Item(43).Delete()
